$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.352.34'
$ws.Range('E2').Value = '  -0.12%  '
$ws.Range('D3').Value = '1.881.29'
$ws.Range('E3').Value = '  +0.32%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7116'
$ws.Range('E5').Value = '  -0.13%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '242.79'
$ws.Range('E6').Value = '  +0.24%  '
$ws.Range('E7').Value = '  +0.13%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.08031'
$ws.Range('E8').Value = '  +2.85%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3140'
$ws.Range('E9').Value = '  +0.74%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '25.11'
$ws.Range('E10').Value = '  -0.31%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08329'
$ws.Range('E11').Value = '  -1.38%  '
$ws.Range('D12').Value = '1.863.34'
$ws.Range('E12').Value = '  -0.68%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.262'
$ws.Range('E13').Value = '  +0.30%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '94.81'
$ws.Range('E14').Value = '  +3.87%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.351'
$ws.Range('E16').Value = '  +4.95%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008690'
$ws.Range('E17').Value = '  +5.28%  '
$ws.Range('D18').Value = '29.373.00'
$ws.Range('E18').Value = '  -0.05%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '243.39'
$ws.Range('E19').Value = '  +0.89%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.36'
$ws.Range('E20').Value = '  +0.48%  '
$ws.Range('D21').Value = '2.147.98'
$ws.Range('E21').Value = '  +1.67%  '
$ws.Range('E22').Value = '  +0.15%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.843'
$ws.Range('E23').Value = '  +0.64%  '
$ws.Range('E25').Value = '  -1.73%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '163.53'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.095'
$ws.Range('E27').Value = '  +0.22%  '
$ws.Range('E28').Value = '  +0.62%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.512'
$ws.Range('E29').Value = '  -0.03%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.439'
$ws.Range('E30').Value = '  +0.12%  '
$ws.Range('E31').Value = '  +0.93%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.206'
$ws.Range('E32').Value = '  -6.45%  '
$ws.Range('E33').Value = '  +1.95%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.942'
$ws.Range('E34').Value = '  -0.06%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7808'
$ws.Range('E35').Value = '  +4.50%  '
$ws.Range('E36').Value = '  -0.05%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.690'
$ws.Range('E37').Value = '  -0.20%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01889'
$ws.Range('E38').Value = '  +0.80%  '
$ws.Range('D39').Value = '1.273.59'
$ws.Range('E39').Value = '  +4.46%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.749'
$ws.Range('E40').Value = '  +0.95%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.548'
$ws.Range('E41').Value = '  +1.11%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9210'
$ws.Range('E42').Value = '  +3.75%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '113.69'
$ws.Range('E43').Value = '  +3.64%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '74.51'
$ws.Range('E44').Value = '  +2.07%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000128'
$ws.Range('E46').Value = '  +4.27%  '
$ws.Range('B47').Value = 'RocketPoolETH'
$ws.Range('C47').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D47').Value = '2.037.89'
$ws.Range('E47').Value = '  +1.01%  '
$ws.Range('E48').Value = '  -0.54%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.5222'
$ws.Range('E49').Value = '  +0.23%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.561'
$ws.Range('E50').Value = '  +1.85%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4384'
$ws.Range('E51').Value = '  +1.30%  '
